$d = $word.ActiveDocument

# Florent : 62h -> Florent : 72h
$d.Content.Find.Execute("62h", $false, $false, $false, $false, $false, $true, 1, $false, "72h", 2) | Out-Null

# Alexis : 58h -> Alexis : 65h
$d.Content.Find.Execute("58h", $false, $false, $false, $false, $false, $true, 1, $false, "65h", 2) | Out-Null
